$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 6
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2

# Remove row 5 entirely (it was deleted in the diff, shrinking the dimension to A1:B4)
$ws.Rows.Item(5).Delete()
